# Updates symbol list data (prices, volumes, and "Hora" hour markers)
# for the cryptos worksheet, per commit "Updated symbol list on
# Mon Jan 30 03:11:24 UTC 2023 with GitHub Actions".
#
# Rows 7-17 also have their Coin name / Link swapped to reflect a
# re-ordering of entries in the source data feed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '315.84' },
    @{ Cell = 'E2'; Value = '2.37%' },
    @{ Cell = 'G2'; Value = '3' },
    @{ Cell = 'D3'; Value = '39.37' },
    @{ Cell = 'E3'; Value = '-1.46%' },
    @{ Cell = 'G3'; Value = '3' },
    @{ Cell = 'D4'; Value = '5.130' },
    @{ Cell = 'E4'; Value = '-0.05%' },
    @{ Cell = 'G4'; Value = '3' },
    @{ Cell = 'D5'; Value = '0.08166' },
    @{ Cell = 'E5'; Value = '0.77%' },
    @{ Cell = 'G5'; Value = '3' },
    @{ Cell = 'D6'; Value = '1.968' },
    @{ Cell = 'E6'; Value = '1.25%' },
    @{ Cell = 'G6'; Value = '3' },
    @{ Cell = 'B7'; Value = 'KuCoinToken' },
    @{ Cell = 'C7'; Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs' },
    @{ Cell = 'D7'; Value = '8.336' },
    @{ Cell = 'E7'; Value = '2.06%' },
    @{ Cell = 'G7'; Value = '3' },
    @{ Cell = 'B8'; Value = 'MXToken' },
    @{ Cell = 'C8'; Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx' },
    @{ Cell = 'D8'; Value = '0.9375' },
    @{ Cell = 'E8'; Value = '1.03%' },
    @{ Cell = 'G8'; Value = '3' },
    @{ Cell = 'B9'; Value = 'LiechtensteinCryptoassetsExchange' },
    @{ Cell = 'C9'; Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx' },
    @{ Cell = 'D9'; Value = '0.1301' },
    @{ Cell = 'E9'; Value = '-8.47%' },
    @{ Cell = 'G9'; Value = '3' },
    @{ Cell = 'B10'; Value = 'WazirX' },
    @{ Cell = 'C10'; Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx' },
    @{ Cell = 'D10'; Value = '0.1969' },
    @{ Cell = 'E10'; Value = '3.04%' },
    @{ Cell = 'G10'; Value = '3' },
    @{ Cell = 'B11'; Value = 'MandalaExchangeToken' },
    @{ Cell = 'C11'; Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx' },
    @{ Cell = 'D11'; Value = '0.09056' },
    @{ Cell = 'E11'; Value = '0.39%' },
    @{ Cell = 'G11'; Value = '3' },
    @{ Cell = 'B12'; Value = 'BitrueCoin' },
    @{ Cell = 'C12'; Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr' },
    @{ Cell = 'D12'; Value = '0.03498' },
    @{ Cell = 'E12'; Value = '-0.82%' },
    @{ Cell = 'G12'; Value = '3' },
    @{ Cell = 'B13'; Value = 'BitMartToken' },
    @{ Cell = 'C13'; Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx' },
    @{ Cell = 'D13'; Value = '0.09758' },
    @{ Cell = 'E13'; Value = '-0.63%' },
    @{ Cell = 'G13'; Value = '3' },
    @{ Cell = 'B14'; Value = 'BitForexToken' },
    @{ Cell = 'C14'; Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf' },
    @{ Cell = 'D14'; Value = '0.001421' },
    @{ Cell = 'E14'; Value = '1.94%' },
    @{ Cell = 'G14'; Value = '3' },
    @{ Cell = 'B15'; Value = 'TigerCash' },
    @{ Cell = 'C15'; Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch' },
    @{ Cell = 'D15'; Value = '0.005998' },
    @{ Cell = 'E15'; Value = '1.60%' },
    @{ Cell = 'G15'; Value = '3' },
    @{ Cell = 'B16'; Value = 'LEO' },
    @{ Cell = 'C16'; Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo' },
    @{ Cell = 'D16'; Value = '3.649' },
    @{ Cell = 'E16'; Value = '-7.65%' },
    @{ Cell = 'G16'; Value = '3' },
    @{ Cell = 'B17'; Value = 'GateToken' },
    @{ Cell = 'C17'; Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt' },
    @{ Cell = 'D17'; Value = '4.376' },
    @{ Cell = 'E17'; Value = '3.45%' },
    @{ Cell = 'G17'; Value = '3' },
    @{ Cell = 'D18'; Value = '3.300' },
    @{ Cell = 'E18'; Value = '-2.73%' },
    @{ Cell = 'G18'; Value = '3' },
    @{ Cell = 'E19'; Value = '1.84%' },
    @{ Cell = 'G19'; Value = '3' },
    @{ Cell = 'D20'; Value = '0.1317' },
    @{ Cell = 'E20'; Value = '-2.34%' },
    @{ Cell = 'G20'; Value = '3' },
    @{ Cell = 'D21'; Value = '4.964' },
    @{ Cell = 'E21'; Value = '6.42%' },
    @{ Cell = 'G21'; Value = '3' },
    @{ Cell = 'D22'; Value = '0.2492' },
    @{ Cell = 'E22'; Value = '-0.31%' },
    @{ Cell = 'G22'; Value = '3' },
    @{ Cell = 'D23'; Value = '0.04356' },
    @{ Cell = 'E23'; Value = '-0.46%' },
    @{ Cell = 'G23'; Value = '3' },
    @{ Cell = 'D24'; Value = '0.001244' },
    @{ Cell = 'E24'; Value = '2.28%' },
    @{ Cell = 'G24'; Value = '3' },
    @{ Cell = 'D25'; Value = '0.004761' },
    @{ Cell = 'E25'; Value = '9.77%' },
    @{ Cell = 'G25'; Value = '3' },
    @{ Cell = 'E26'; Value = '199.30%' },
    @{ Cell = 'G26'; Value = '3' },
    @{ Cell = 'E27'; Value = '-7.71%' },
    @{ Cell = 'G27'; Value = '3' },
    @{ Cell = 'G28'; Value = '3' },
    @{ Cell = 'G29'; Value = '3' },
    @{ Cell = 'G30'; Value = '3' },
    @{ Cell = 'G31'; Value = '3' },
    @{ Cell = 'G32'; Value = '3' },
    @{ Cell = 'G33'; Value = '3' },
    @{ Cell = 'G34'; Value = '3' },
    @{ Cell = 'G35'; Value = '3' },
    @{ Cell = 'G36'; Value = '3' },
    @{ Cell = 'G37'; Value = '3' },
    @{ Cell = 'G38'; Value = '3' },
    @{ Cell = 'D39'; Value = '0.02196' },
    @{ Cell = 'E39'; Value = '8.47%' },
    @{ Cell = 'G39'; Value = '3' },
    @{ Cell = 'D40'; Value = '0.05175' },
    @{ Cell = 'E40'; Value = '2.53%' },
    @{ Cell = 'G40'; Value = '3' },
    @{ Cell = 'D41'; Value = '0.007767' },
    @{ Cell = 'E41'; Value = '5.04%' },
    @{ Cell = 'G41'; Value = '3' },
    @{ Cell = 'D42'; Value = '0.01040' },
    @{ Cell = 'E42'; Value = '6.40%' },
    @{ Cell = 'G42'; Value = '3' },
    @{ Cell = 'D43'; Value = '0.1401' },
    @{ Cell = 'E43'; Value = '2.66%' },
    @{ Cell = 'G43'; Value = '3' },
    @{ Cell = 'D44'; Value = '0.002083' },
    @{ Cell = 'E44'; Value = '-2.32%' },
    @{ Cell = 'G44'; Value = '3' },
    @{ Cell = 'D45'; Value = '0.009289' },
    @{ Cell = 'E45'; Value = '1.37%' },
    @{ Cell = 'G45'; Value = '3' },
    @{ Cell = 'E46'; Value = '8.89%' },
    @{ Cell = 'G46'; Value = '3' },
    @{ Cell = 'D47'; Value = '0.00000000751' },
    @{ Cell = 'E47'; Value = '-0.03%' },
    @{ Cell = 'G47'; Value = '3' },
    @{ Cell = 'D48'; Value = '0.002885' },
    @{ Cell = 'E48'; Value = '0.73%' },
    @{ Cell = 'G48'; Value = '3' },
    @{ Cell = 'D49'; Value = '0.001693' },
    @{ Cell = 'E49'; Value = '30.16%' },
    @{ Cell = 'G49'; Value = '3' },
    @{ Cell = 'D50'; Value = '0.00002103' },
    @{ Cell = 'E50'; Value = '-0.03%' },
    @{ Cell = 'G50'; Value = '3' },
    @{ Cell = 'D51'; Value = '0.0002003' },
    @{ Cell = 'E51'; Value = '-0.03%' },
    @{ Cell = 'G51'; Value = '3' }
)

foreach ($item in $updates) {
    $cell = $ws.Range($item.Cell)
    # Force a text number format before assigning so that Excel does not
    # reinterpret numeric-looking / percent-looking strings (e.g. "315.84",
    # "2.37%", "3") as actual numbers and mangle their literal formatting
    # (trailing zeros, percent scaling, etc.).
    $cell.NumberFormat = "@"
    $cell.Value = $item.Value
    # Reset back to the workbook's default "Normal" style so we don't leave
    # a stray text-format style applied to the cell (matches original file,
    # where these cells carry no explicit style index).
    $cell.Style = "Normal"
}
